# "Werkelijk resultaat toegevoegd aan werkbestand zodat syncen nu goed gaat"
#
# Sheet1's header row (row 2) had "Verwacht resultaat" in G2 and
# "Aangepast resultaat" in H2, which put the new "Werkelijk resultaat"
# column (F2) out of sync with the rest of the automation. Swap the two
# headers so the column order lines up correctly again: "Aangepast
# resultaat" becomes G2 and "Verwacht resultaat" becomes H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = "Aangepast resultaat"
$ws.Range("H2").Value = "Verwacht resultaat"

# F2 ("Werkelijk resultaat") carried a one-off bordered style left over from
# earlier edits; drop the border so it matches the plain header formatting
# used by the other non-bordered header cells.
$ws.Range("F2").Borders.LineStyle = -4142

# Column G's content got longer ("Aangepast resultaat"), so re-fit its
# width to the new text, same as column H.
$ws.Range("G1:G2").EntireColumn.AutoFit()

# Park the selection near the top of the sheet.
$ws.Range("H3").Select()
